# Apply the commit's changes to the betexplorer Turkey 1-Lig workbook.
#
# 1) A handful of existing row-pairs had their match data (columns B:V)
#    swapped between the two rows (column A / Indice stays put).
# 2) Nine brand-new match rows (146-154) are appended at the bottom,
#    copying the formatting of row 145 for column A (bold/centered index)
#    and column E (date/time number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($sheet, $rowA, $rowB) {
    for ($col = 2; $col -le 22; $col++) {
        $cellA = $sheet.Cells.Item($rowA, $col)
        $cellB = $sheet.Cells.Item($rowB, $col)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# --- 1) Swap the row-pairs whose home/away data got exchanged ---
Swap-MatchRows $ws 15 16
Swap-MatchRows $ws 36 37
Swap-MatchRows $ws 42 43
Swap-MatchRows $ws 44 45
Swap-MatchRows $ws 47 48
Swap-MatchRows $ws 53 54
Swap-MatchRows $ws 79 80
Swap-MatchRows $ws 84 85
Swap-MatchRows $ws 89 90
Swap-MatchRows $ws 96 97

# --- 2) Append the new rows (146-154) ---

function Add-MatchRow($sheet, $row, $indice, $dataPartida, $home, $homeGols, $away, $awayGols,
    $homeOpenOdds, $homeOpenDt, $homeCloseOdds, $homeCloseDt,
    $drawOpenOdds, $drawOpenDt, $drawCloseOdds, $drawCloseDt,
    $awayOpenOdds, $awayOpenDt, $awayCloseOdds, $awayCloseDt, $url) {

    # Copy formatting from the last pre-existing row (145) so the new
    # rows keep the same styles (bold/centered index, date number format).
    $sheet.Range("A145").Copy()
    $sheet.Range("A$row").PasteSpecial(-4122)
    $sheet.Range("E145").Copy()
    $sheet.Range("E$row").PasteSpecial(-4122)

    $sheet.Cells.Item($row, 1).Value = $indice
    $sheet.Cells.Item($row, 2).Value = "turkey"
    $sheet.Cells.Item($row, 3).Value = "1-lig"
    $sheet.Cells.Item($row, 4).Value = "2023-2024"
    $sheet.Cells.Item($row, 5).Value = $dataPartida
    $sheet.Cells.Item($row, 6).Value = $home
    $sheet.Cells.Item($row, 7).Value = $homeGols
    $sheet.Cells.Item($row, 8).Value = $away
    $sheet.Cells.Item($row, 9).Value = $awayGols
    $sheet.Cells.Item($row, 10).Value = $homeOpenOdds
    $sheet.Cells.Item($row, 11).Value = $homeOpenDt
    $sheet.Cells.Item($row, 12).Value = $homeCloseOdds
    $sheet.Cells.Item($row, 13).Value = $homeCloseDt
    $sheet.Cells.Item($row, 14).Value = $drawOpenOdds
    $sheet.Cells.Item($row, 15).Value = $drawOpenDt
    $sheet.Cells.Item($row, 16).Value = $drawCloseOdds
    $sheet.Cells.Item($row, 17).Value = $drawCloseDt
    $sheet.Cells.Item($row, 18).Value = $awayOpenOdds
    $sheet.Cells.Item($row, 19).Value = $awayOpenDt
    $sheet.Cells.Item($row, 20).Value = $awayCloseOdds
    $sheet.Cells.Item($row, 21).Value = $awayCloseDt
    $sheet.Cells.Item($row, 22).Value = $url
}

Add-MatchRow $ws 146 145 45283.47916666666 "Genclerbirligi" 0 "Tuzlaspor" 2 `
    1.61 "20/12/2023 09:42" 1.71 "23/12/2023 11:26" `
    3.79 "20/12/2023 09:42" 3.61 "23/12/2023 11:29" `
    5.01 "20/12/2023 09:42" 5.32 "23/12/2023 11:29" `
    "https://www.betexplorer.com/football/turkey/1-lig/genclerbirligi-tuzlaspor/6ZKcRN2t/"

Add-MatchRow $ws 147 146 45283.70833333334 "Sakaryaspor" 3 "Kocaelispor" 1 `
    2.34 "19/12/2023 18:13" 2.59 "23/12/2023 16:56" `
    3.17 "19/12/2023 18:13" 3.36 "23/12/2023 16:55" `
    3.05 "19/12/2023 18:13" 2.78 "23/12/2023 16:56" `
    "https://www.betexplorer.com/football/turkey/1-lig/sakaryaspor-kocaelispor/lbA2QsIn/"

Add-MatchRow $ws 148 147 45284.47916666666 "Bandirmaspor" 2 "Adanaspor AS" 3 `
    1.41 "20/12/2023 15:12" 1.38 "24/12/2023 11:28" `
    4.48 "20/12/2023 15:12" 4.84 "24/12/2023 11:29" `
    6.64 "20/12/2023 15:12" 8.34 "24/12/2023 11:29" `
    "https://www.betexplorer.com/football/turkey/1-lig/bandirmaspor-adanaspor-as/Iw2JMaHB/"

Add-MatchRow $ws 149 148 45284.47916666666 "Bodrumspor" 3 "Goztepe" 0 `
    2.1 "20/12/2023 18:12" 2.93 "24/12/2023 11:29" `
    3.05 "20/12/2023 18:12" 3.23 "24/12/2023 11:29" `
    3.72 "20/12/2023 18:12" 2.55 "24/12/2023 11:29" `
    "https://www.betexplorer.com/football/turkey/1-lig/bodrumspor-goztepe/dODAOLna/"

Add-MatchRow $ws 150 149 45284.58333333334 "Umraniyespor" 2 "Giresunspor" 1 `
    1.45 "20/12/2023 12:12" 1.28 "24/12/2023 13:33" `
    4.28 "20/12/2023 12:12" 5.82 "24/12/2023 13:40" `
    6.18 "20/12/2023 12:12" 10.15 "24/12/2023 13:40" `
    "https://www.betexplorer.com/football/turkey/1-lig/umraniyespor-giresunspor/Sl96P1Xh/"

Add-MatchRow $ws 151 150 45284.70833333334 "Altay" 2 "Keciorengucu" 1 `
    3.68 "20/12/2023 15:12" 4.08 "24/12/2023 16:57" `
    3.45 "20/12/2023 15:12" 3.79 "24/12/2023 16:59" `
    1.93 "20/12/2023 15:12" 1.85 "24/12/2023 16:57" `
    "https://www.betexplorer.com/football/turkey/1-lig/altay-keciorengucu/OxUj4i1p/"

Add-MatchRow $ws 152 151 45285.5 "Manisa FK" 1 "Sanliurfaspor" 1 `
    1.49 "21/12/2023 15:12" 1.49 "25/12/2023 11:58" `
    4.16 "21/12/2023 15:12" 4.22 "25/12/2023 11:58" `
    5.79 "21/12/2023 15:12" 6.98 "25/12/2023 11:58" `
    "https://www.betexplorer.com/football/turkey/1-lig/manisa-fk-sanliurfaspor/2TSf3BGj/"

Add-MatchRow $ws 153 152 45285.625 "Corum" 2 "Boluspor" 0 `
    1.98 "21/12/2023 15:12" 2 "25/12/2023 14:55" `
    3.31 "21/12/2023 15:12" 3.4 "25/12/2023 14:55" `
    3.73 "21/12/2023 15:12" 3.96 "25/12/2023 14:55" `
    "https://www.betexplorer.com/football/turkey/1-lig/corum-fk-boluspor/zFCENu25/"

Add-MatchRow $ws 154 153 45285.75 "Eyupspor" 3 "Erzurumspor" 0 `
    1.25 "20/12/2023 18:12" 1.33 "25/12/2023 17:52" `
    5.28 "20/12/2023 18:12" 5.39 "25/12/2023 17:57" `
    9.93 "20/12/2023 18:12" 8.84 "25/12/2023 17:57" `
    "https://www.betexplorer.com/football/turkey/1-lig/eyupspor-erzurumspor-fk/Wt6NLJWH/"
